# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.431.10'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '3.557.55'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.67'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.35'
$ws.Range('E6').Value = '  +4.01%  '
$ws.Range('D7').Value = '3.550.21'
$ws.Range('E7').Value = '  +2.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.620'
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +11.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.645'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.16'
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000311'
$ws.Range('E13').Value = '  +3.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.40'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '4.125.13'
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('D16').Value = '70.460.94'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.585.94'
$ws.Range('E17').Value = '  +2.93%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.76'
$ws.Range('E18').Value = '  +5.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.98'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '582.24'
$ws.Range('E20').Value = '  +9.25%  '
$ws.Range('E21').Value = '  +1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.993'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.93'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('E24').Value = '  +4.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.85'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.52'
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.00'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.92'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.35'
$ws.Range('E29').Value = '  +4.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.31'
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.05'
$ws.Range('E31').Value = '  -1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.19'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.36'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('B35').Value = 'dogwifhat'
$ws.Range('C35').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.80'
$ws.Range('E35').Value = '  +29.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.22'
$ws.Range('E36').Value = '  +6.85%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '527.02'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.406'
$ws.Range('E38').Value = '  +1.54%  '
$ws.Range('D39').Value = '3.685.21'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.92'
$ws.Range('E40').Value = '  +1.15%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '0.0₃0788'
$ws.Range('E42').Value = '  +6.35%  '
$ws.Range('E43').Value = '  +5.96%  '
$ws.Range('E44').Value = '  +3.55%  '
$ws.Range('E45').Value = '  +5.69%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.94'
$ws.Range('E46').Value = '  +0.89%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('E48').Value = '  +3.75%  '
$ws.Range('E49').Value = '  +3.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.43'
$ws.Range('E51').Value = '  +9.21%  '
